$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Activate()

# 1) Update the pidRS(<int>) command doc to use <rs> and mention labeled patterns.
$ws.Range("B80").Value = "pidRS(<rs>)"
$ws.Range("C80").Value = "activates the PID Ramp-Soak pattern number <rs> (1-based!) or the one labeled <rs>"

# 2) Update the palette(<int>) command doc to use <p> and mention palette labels.
$ws.Range("B87").Value = "palette(<p>)"
$ws.Range("C87").Value = "activates palette <p> with <p> either a number 0-9 or a palette label"

# 3) Insert a new row right after "clearBackground" (row 91) for the new
#    alarmset(<as>) command, pushing the RC Command block (and everything
#    below it) down by one row.
$ws.Rows.Item(92).Insert()
$ws.Range("B92").Value = "alarmset(<as>)"
$ws.Range("C92").Value = "activates the alarmset with the given number or label"

# Widen column C a bit to better fit the new/longer descriptions.
$ws.Columns.Item(3).ColumnWidth = 26.15

# Leave the selection where the author left it after adding the new row.
$ws.Range("C87").Select()
